# Banco_Dashboard.xlsx - progress update on "dados_corrigidos" sheet.
# Mirrors a manual edit session in Excel: a handful of activities had their
# M2_Realizado (column D) - and for one of them also M2_Previsto (column C) -
# bumped up, which flips their Status (column E, a shared formula) from
# "Não iniciado" to "Em andamento". The analyst also widened column B to read
# the full activity names and left the cursor on D8 when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - M2_Previsto increased as well as M2_Realizado.
$ws.Range("C7").Value = 120
$ws.Range("D7").Value = 60

# Row 16 - only M2_Realizado changed.
$ws.Range("D16").Value = 5.7

# Row 20
$ws.Range("D20").Value = 18

# Row 63
$ws.Range("D63").Value = 1

# Row 71
$ws.Range("D71").Value = 12

# Widen column B (Atividade) so the longer activity descriptions fit.
$ws.Columns.Item(2).ColumnWidth = 92.21875

# Leave the active selection where the analyst left it before saving.
$ws.Range("D8").Select() | Out-Null
